# Generate Report for Handoff
# Adds two new localization-status rows (for files
# 6b25841a-d90f-441c-905d-c757e590abdc.md and
# fc55cf56-39b8-40b2-9d02-79d455db6881.md) to the Overview / zh-cn / de-de
# sheets of the localization-status workbook: one inserted right after the
# existing "eedced51" row, one appended at the end. The pre-existing
# "84ac6ebe" row shifts down by one row in the process.

$wb = $excel.ActiveWorkbook

$dateFmt = "yyyy-mm-dd HH:mm:ss"

function Set-Cell($sheet, $row, $col, $value) {
    $sheet.Cells.Item($row, $col).Value = $value
}

# ---------------------------------------------------------------------
# Sheet "Overview": columns A..G
#   A File Name   B Path And Name (hyperlink)  C Extension  D Publish URL
#   E zh-cn       F de-de                       G Latest HO Xliff Generate Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Row 3 becomes the new "6b25841a" entry (was "84ac6ebe").
Set-Cell $wsOverview 3 1 "6b25841a-d90f-441c-905d-c757e590abdc.md"
Set-Cell $wsOverview 3 2 "e2e\6b25841a-d90f-441c-905d-c757e590abdc.md"
Set-Cell $wsOverview 3 3 ".md"
Set-Cell $wsOverview 3 5 "Ready for handoff"
Set-Cell $wsOverview 3 6 "Ready for handoff"
Set-Cell $wsOverview 3 7 "2016-08-29 14:46:36"
$wsOverview.Cells.Item(3, 7).NumberFormat = $dateFmt

# Row 4: the "84ac6ebe" entry, shifted down from the old row 3.
Set-Cell $wsOverview 4 1 "84ac6ebe-2963-4179-802c-a62a53aaa5f2.md"
Set-Cell $wsOverview 4 2 "e2e\84ac6ebe-2963-4179-802c-a62a53aaa5f2.md"
Set-Cell $wsOverview 4 3 ".md"
Set-Cell $wsOverview 4 5 "Ready for handoff"
Set-Cell $wsOverview 4 6 "Ready for handoff"
Set-Cell $wsOverview 4 7 "2016-08-29 14:44:53"
$wsOverview.Cells.Item(4, 7).NumberFormat = $dateFmt

# Row 5: new "fc55cf56" entry, appended.
Set-Cell $wsOverview 5 1 "fc55cf56-39b8-40b2-9d02-79d455db6881.md"
Set-Cell $wsOverview 5 2 "e2e\fc55cf56-39b8-40b2-9d02-79d455db6881.md"
Set-Cell $wsOverview 5 3 ".md"
Set-Cell $wsOverview 5 5 "Ready for handoff"
Set-Cell $wsOverview 5 6 "Ready for handoff"
Set-Cell $wsOverview 5 7 "2016-08-29 14:46:36"
$wsOverview.Cells.Item(5, 7).NumberFormat = $dateFmt

# Hyperlinks on column B for rows 3..5 (re-create the ones that moved,
# and add the brand new ones).
$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item(3, 2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68c5439b8ea4cacc3f45e9f3f278e98c2c1bd750/e2e/6b25841a-d90f-441c-905d-c757e590abdc.md", "", "", "e2e\6b25841a-d90f-441c-905d-c757e590abdc.md")
$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item(4, 2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2f71f7d7c8356ede91e9ca8d46b178724adabf91/e2e/84ac6ebe-2963-4179-802c-a62a53aaa5f2.md", "", "", "e2e\84ac6ebe-2963-4179-802c-a62a53aaa5f2.md")
$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item(5, 2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4f12d58c0fc3502ad6454aa803df1acaf4ca16c5/e2e/fc55cf56-39b8-40b2-9d02-79d455db6881.md", "", "", "e2e\fc55cf56-39b8-40b2-9d02-79d455db6881.md")

# Grow the "Overview" table + used dimension to A1:G5.
$loOverview = $wsOverview.ListObjects.Item("Overview")
$loOverview.Resize($wsOverview.Range("A1:G5"))

# ---------------------------------------------------------------------
# Sheets "zh-cn" / "de-de": columns A..P
#   A Source File Name (hyperlink)   B File Extension   C Status
#   D Source Path   E Priority   F Content Duplicate
#   G Latest Handoff File   H Latest Handoff Datetime   I Latest Target File
#   J Latest Handback File  K Latest Handback DateTime  L Reference Tokens
#   M To be localized       N Dependency From           O Has metadata
#   P Error Detail
# ---------------------------------------------------------------------
function Update-LangSheet($sheetName, $xlfSuffix, $row3Handoff, $row3HandoffDate, $row4Handoff, $row4HandoffDate, $row5Handoff, $row5HandoffDate, $hlShaNew, $hlShaExisting, $hlShaAppended) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 3: new "6b25841a" entry.
    Set-Cell $ws 3 1 "6b25841a-d90f-441c-905d-c757e590abdc.md"
    Set-Cell $ws 3 2 ".md"
    Set-Cell $ws 3 3 "Ready for handoff"
    Set-Cell $ws 3 4 "e2e"
    Set-Cell $ws 3 5 "ht"
    Set-Cell $ws 3 6 "False"
    Set-Cell $ws 3 7 $row3Handoff
    Set-Cell $ws 3 8 $row3HandoffDate
    $ws.Cells.Item(3, 8).NumberFormat = $dateFmt
    Set-Cell $ws 3 11 "0001-01-01 00:00:00"
    $ws.Cells.Item(3, 11).NumberFormat = $dateFmt
    Set-Cell $ws 3 13 "True"
    Set-Cell $ws 3 15 "False"

    # Row 4: "84ac6ebe" entry, shifted down from the old row 3.
    Set-Cell $ws 4 1 "84ac6ebe-2963-4179-802c-a62a53aaa5f2.md"
    Set-Cell $ws 4 2 ".md"
    Set-Cell $ws 4 3 "Ready for handoff"
    Set-Cell $ws 4 4 "e2e"
    Set-Cell $ws 4 5 "ht"
    Set-Cell $ws 4 6 "False"
    Set-Cell $ws 4 7 $row4Handoff
    Set-Cell $ws 4 8 $row4HandoffDate
    $ws.Cells.Item(4, 8).NumberFormat = $dateFmt
    Set-Cell $ws 4 11 "0001-01-01 00:00:00"
    $ws.Cells.Item(4, 11).NumberFormat = $dateFmt
    Set-Cell $ws 4 13 "True"
    Set-Cell $ws 4 15 "False"

    # Row 5: new "fc55cf56" entry, appended.
    Set-Cell $ws 5 1 "fc55cf56-39b8-40b2-9d02-79d455db6881.md"
    Set-Cell $ws 5 2 ".md"
    Set-Cell $ws 5 3 "Ready for handoff"
    Set-Cell $ws 5 4 "e2e"
    Set-Cell $ws 5 5 "ht"
    Set-Cell $ws 5 6 "False"
    Set-Cell $ws 5 7 $row5Handoff
    Set-Cell $ws 5 8 $row5HandoffDate
    $ws.Cells.Item(5, 8).NumberFormat = $dateFmt
    Set-Cell $ws 5 11 "0001-01-01 00:00:00"
    $ws.Cells.Item(5, 11).NumberFormat = $dateFmt
    Set-Cell $ws 5 13 "True"
    Set-Cell $ws 5 15 "False"

    # Hyperlinks on column A for rows 3..5.
    $ws.Hyperlinks.Add($ws.Cells.Item(3, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68c5439b8ea4cacc3f45e9f3f278e98c2c1bd750/e2e/6b25841a-d90f-441c-905d-c757e590abdc.md", "", "", "6b25841a-d90f-441c-905d-c757e590abdc.md")
    $ws.Hyperlinks.Add($ws.Cells.Item(4, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2f71f7d7c8356ede91e9ca8d46b178724adabf91/e2e/84ac6ebe-2963-4179-802c-a62a53aaa5f2.md", "", "", "84ac6ebe-2963-4179-802c-a62a53aaa5f2.md")
    $ws.Hyperlinks.Add($ws.Cells.Item(5, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4f12d58c0fc3502ad6454aa803df1acaf4ca16c5/e2e/fc55cf56-39b8-40b2-9d02-79d455db6881.md", "", "", "fc55cf56-39b8-40b2-9d02-79d455db6881.md")

    # Grow the language table + used dimension to A1:P5.
    $lo = $ws.ListObjects.Item($sheetName)
    $lo.Resize($ws.Range("A1:P5"))
}

Update-LangSheet "zh-cn" "zh-cn.xlf" `
    "6b25841a-d90f-441c-905d-c757e590abdc.68c5439b8ea4cacc3f45e9f3f278e98c2c1bd750.zh-cn.xlf" "2016-08-29 14:46:31" `
    "84ac6ebe-2963-4179-802c-a62a53aaa5f2.00917b7056fea0cefe0d69a71198b7df9afc3a23.zh-cn.xlf" "2016-08-29 14:44:49" `
    "fc55cf56-39b8-40b2-9d02-79d455db6881.4f12d58c0fc3502ad6454aa803df1acaf4ca16c5.zh-cn.xlf" "2016-08-29 14:46:31" `
    "" "" ""

Update-LangSheet "de-de" "de-de.xlf" `
    "6b25841a-d90f-441c-905d-c757e590abdc.68c5439b8ea4cacc3f45e9f3f278e98c2c1bd750.de-de.xlf" "2016-08-29 14:46:36" `
    "84ac6ebe-2963-4179-802c-a62a53aaa5f2.00917b7056fea0cefe0d69a71198b7df9afc3a23.de-de.xlf" "2016-08-29 14:44:53" `
    "fc55cf56-39b8-40b2-9d02-79d455db6881.4f12d58c0fc3502ad6454aa803df1acaf4ca16c5.de-de.xlf" "2016-08-29 14:46:36" `
    "" "" ""

Write-Output "Localization status report updated: 2 new rows added (6b25841a-d90f-441c-905d-c757e590abdc, fc55cf56-39b8-40b2-9d02-79d455db6881)."
